{"js": "// Replace each exact two-digit division answer string in the document's\n// table cells with its new value. Every \"find\" string below occurs exactly\n// once in the document (verified against before.docx), so a body-wide\n// search-and-replace keyed on the old text is unambiguous and targets the\n// correct cell without relying on row/column indices.\n//\n// The pairs are applied in the same order they appear in the source\n// document. That matters because one pair's replacement text (\"62\u00f73=20, 2\")\n// is identical to an earlier pair's find text, so processing in document\n// order guarantees each search still resolves to a single, correct hit.\nconst replacements = [\n  [\"18\u00f74=4, 2\", \"93\u00f78=11, 5\"],\n  [\"62\u00f73=20, 2\", \"76\u00f79=8, 4\"],\n  [\"85\u00f78=10, 5\", \"12\u00f73=4, 0\"],\n  [\"80\u00f76=13, 2\", \"73\u00f75=14, 3\"],\n  [\"37\u00f73=12, 1\", \"92\u00f72=46, 0\"],\n  [\"98\u00f76=16, 2\", \"99\u00f77=14, 1\"],\n  [\"85\u00f73=28, 1\", \"68\u00f75=13, 3\"],\n  [\"75\u00f74=18, 3\", \"38\u00f72=19, 0\"],\n  [\"82\u00f73=27, 1\", \"21\u00f75=4, 1\"],\n  [\"69\u00f73=23, 0\", \"75\u00f73=25, 0\"],\n  [\"15\u00f72=7, 1\", \"44\u00f76=7, 2\"],\n  [\"52\u00f74=13, 0\", \"41\u00f72=20, 1\"],\n  [\"21\u00f74=5, 1\", \"59\u00f74=14, 3\"],\n  [\"74\u00f72=37, 0\", \"28\u00f77=4, 0\"],\n  [\"81\u00f78=10, 1\", \"62\u00f73=20, 2\"],\n  [\"39\u00f77=5, 4\", \"24\u00f77=3, 3\"],\n  [\"13\u00f76=2, 1\", \"57\u00f76=9, 3\"],\n  [\"23\u00f79=2, 5\", \"39\u00f72=19, 1\"],\n  [\"56\u00f74=14, 0\", \"88\u00f77=12, 4\"],\n  [\"94\u00f79=10, 4\", \"84\u00f74=21, 0\"],\n  [\"17\u00f76=2, 5\", \"64\u00f76=10, 4\"],\n  [\"69\u00f76=11, 3\", \"37\u00f78=4, 5\"],\n  [\"93\u00f75=18, 3\", \"22\u00f77=3, 1\"],\n  [\"35\u00f75=7, 0\", \"92\u00f74=23, 0\"],\n  [\"28\u00f72=14, 0\", \"29\u00f74=7, 1\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(`Expected exactly one match for \"${oldText}\", found ${results.items.length}`);\n  }\n\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Replace each exact two-digit division answer string in the document's\n# table cells with its new value. Every \"find\" string occurs exactly once\n# in the document (verified against before.docx), so Find/Replace keyed on\n# the old text is unambiguous and needs no row/column bookkeeping. The pairs\n# are applied in the same order as they appear in the source document so that\n# an old value that is itself the replacement text of an earlier pair (e.g.\n# \"62\u00f73=20, 2\") is never mistakenly re-matched.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"18\u00f74=4, 2\", \"93\u00f78=11, 5\"),\n  @(\"62\u00f73=20, 2\", \"76\u00f79=8, 4\"),\n  @(\"85\u00f78=10, 5\", \"12\u00f73=4, 0\"),\n  @(\"80\u00f76=13, 2\", \"73\u00f75=14, 3\"),\n  @(\"37\u00f73=12, 1\", \"92\u00f72=46, 0\"),\n  @(\"98\u00f76=16, 2\", \"99\u00f77=14, 1\"),\n  @(\"85\u00f73=28, 1\", \"68\u00f75=13, 3\"),\n  @(\"75\u00f74=18, 3\", \"38\u00f72=19, 0\"),\n  @(\"82\u00f73=27, 1\", \"21\u00f75=4, 1\"),\n  @(\"69\u00f73=23, 0\", \"75\u00f73=25, 0\"),\n  @(\"15\u00f72=7, 1\", \"44\u00f76=7, 2\"),\n  @(\"52\u00f74=13, 0\", \"41\u00f72=20, 1\"),\n  @(\"21\u00f74=5, 1\", \"59\u00f74=14, 3\"),\n  @(\"74\u00f72=37, 0\", \"28\u00f77=4, 0\"),\n  @(\"81\u00f78=10, 1\", \"62\u00f73=20, 2\"),\n  @(\"39\u00f77=5, 4\", \"24\u00f77=3, 3\"),\n  @(\"13\u00f76=2, 1\", \"57\u00f76=9, 3\"),\n  @(\"23\u00f79=2, 5\", \"39\u00f72=19, 1\"),\n  @(\"56\u00f74=14, 0\", \"88\u00f77=12, 4\"),\n  @(\"94\u00f79=10, 4\", \"84\u00f74=21, 0\"),\n  @(\"17\u00f76=2, 5\", \"64\u00f76=10, 4\"),\n  @(\"69\u00f76=11, 3\", \"37\u00f78=4, 5\"),\n  @(\"93\u00f75=18, 3\", \"22\u00f77=3, 1\"),\n  @(\"35\u00f75=7, 0\", \"92\u00f74=23, 0\"),\n  @(\"28\u00f72=14, 0\", \"29\u00f74=7, 1\")\n)\n\n# Find.Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,\n#              MatchSoundsLike, MatchAllWordForms, Forward, Wrap, Format,\n#              ReplaceWith, Replace) -- wdFindContinue=1, wdReplaceAll=2.\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $find = $d.Content.Find\n  $found = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n  if (-not $found) {\n    throw \"Text not found: $oldText\"\n  }\n}\n\n"}
